$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: date moves from 42060 -> 42061, add hours/role/activity
# (role/activity duplicate row 9's "Test Analyst" / "Test case writing")
$ws.Range("A10").Value = 42061
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Test Analyst"
$ws.Range("D10").Value = "Test case writing"

# Row 11: date moves from 42060 -> 42061, add hours/role/activity
$ws.Range("A11").Value = 42061
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "Requirements Specifier"
$ws.Range("D11").Value = "Aktivitetsdiagram, domænemodel"

# Row 12 was blank; A12 needs the same date format/style as the other date
# cells, so copy formatting from an existing date cell before writing the
# value. Fill D before C so new shared strings are created in the same
# order as the target workbook.
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = 42062
$ws.Range("B12").Value = 2.5
$ws.Range("D12").Value = "Test case implementation"
$ws.Range("C12").Value = "Test Designer"

# Row 13 was blank; same formatting fix as row 12.
$ws.Range("A9").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 42062
$ws.Range("B13").Value = 0.5
$ws.Range("C13").Value = "Any role"
$ws.Range("D13").Value = "GitHub"

# Row 14 was blank; only the date gets filled in (same formatting fix).
$ws.Range("A9").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 42062

$ws.Range("A14").Select() | Out-Null
